# Auto-generated edit script applying numeric updates per the commit diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 71432536
$ws.Range("I113").Value = 200002740
$ws.Range("K113").Value = 200002740
$ws.Range("M113").Value = -199999486

$ws.Range("H130").Value = 69500
$ws.Range("J130").Value = 79000
$ws.Range("L130").Value = 79000
$ws.Range("N130").Value = -89040

$ws.Range("H132").Value = 3878.6272
$ws.Range("I132").Value = 3924.3618
$ws.Range("J132").Value = 3699.5
$ws.Range("K132").Value = 11773.0854
$ws.Range("L132").Value = 11098.5
$ws.Range("M132").Value = -9243.0854
$ws.Range("N132").Value = -16158.5

$ws.Range("H137").Value = 3220.4595
$ws.Range("J137").Value = 4314.35
$ws.Range("L137").Value = 12943.05
$ws.Range("N137").Value = -18043.05

$ws.Range("H138").Value = 2559.4146
$ws.Range("I138").Value = 1862.2593
$ws.Range("J138").Value = 3903.9285
$ws.Range("K138").Value = 5586.7779
$ws.Range("L138").Value = 11711.7855
$ws.Range("M138").Value = -446.7779
$ws.Range("N138").Value = -21991.7855

$ws.Range("H139").Value = 55166.5
$ws.Range("J139").Value = 59199.8
$ws.Range("L139").Value = 59199.8
$ws.Range("N139").Value = -69479.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 190872.38
$ws.Range("I32").Value = 190872.38
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 190872.38
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -190585.38
$ws.Range("N32").ClearContents()

$ws.Range("H45").Value = 372570.53
$ws.Range("I45").Value = 557273.6
$ws.Range("K45").Value = 557273.6
$ws.Range("M45").Value = -556896.6

$ws.Range("H63").Value = 70599840
$ws.Range("I63").Value = 142861900
$ws.Range("K63").Value = 142861900
$ws.Range("M63").Value = -142861214

$ws.Range("H66").Value = 70599840
$ws.Range("I66").Value = 142861900
$ws.Range("K66").Value = 714309500
$ws.Range("M66").Value = -714306068

$ws.Range("H74").Value = 3216.1353
$ws.Range("I74").Value = 3454.2415
$ws.Range("J74").Value = 2353
$ws.Range("K74").Value = 3454.2415
$ws.Range("L74").Value = 2353
$ws.Range("M74").Value = -2580.2415
$ws.Range("N74").Value = -4101

$ws.Range("H77").Value = 3216.1353
$ws.Range("I77").Value = 3454.2415
$ws.Range("J77").Value = 2353
$ws.Range("K77").Value = 17271.2075
$ws.Range("L77").Value = 11765
$ws.Range("M77").Value = -12903.2075
$ws.Range("N77").Value = -20501

$ws.Range("H88").Value = 18521304
$ws.Range("I88").Value = 83334080
$ws.Range("J88").Value = 3367.5715
$ws.Range("K88").Value = 83334080
$ws.Range("L88").Value = 3367.5715
$ws.Range("M88").Value = -83333674
$ws.Range("N88").Value = -4179.5715

$ws.Range("H91").Value = 18521304
$ws.Range("I91").Value = 83334080
$ws.Range("J91").Value = 3367.5715
$ws.Range("K91").Value = 83334080
$ws.Range("L91").Value = 3367.5715
$ws.Range("M91").Value = -83332676
$ws.Range("N91").Value = -6175.5715

$ws.Range("H110").Value = 37043536
$ws.Range("I110").Value = 43479108
$ws.Range("K110").Value = 43479108
$ws.Range("M110").Value = -43477063

$ws.Range("H122").Value = 27781946
$ws.Range("I122").Value = 66669468
$ws.Range("K122").Value = 200008404
$ws.Range("M122").Value = -200005954

$ws.Range("H132").Value = 58825980
$ws.Range("I132").Value = 142859490
$ws.Range("K132").Value = 428578470
$ws.Range("M132").Value = -428575940

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2161
$ws.Range("I86").Value = 2299.625
$ws.Range("J86").Value = 2068.5833
$ws.Range("K86").Value = 2299.625
$ws.Range("L86").Value = 2068.5833
$ws.Range("M86").Value = -1176.625
$ws.Range("N86").Value = -4314.5833

$ws.Range("H89").Value = 2161
$ws.Range("I89").Value = 2299.625
$ws.Range("J89").Value = 2068.5833
$ws.Range("K89").Value = 11498.125
$ws.Range("L89").Value = 10342.9165
$ws.Range("M89").Value = -5882.125
$ws.Range("N89").Value = -21574.9165

$ws.Range("H134").Value = 2364.9092
$ws.Range("I134").Value = 2091.7073
$ws.Range("J134").Value = 6098.6665
$ws.Range("K134").Value = 6275.1219
$ws.Range("L134").Value = 18295.9995
$ws.Range("M134").Value = -3740.1219
$ws.Range("N134").Value = -23365.9995

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3788.5386
$ws.Range("J31").Value = 4025.5667
$ws.Range("L31").Value = 4025.5667
$ws.Range("N31").Value = -4615.566699999999

$ws.Range("H34").Value = 3788.5386
$ws.Range("J34").Value = 4025.5667
$ws.Range("L34").Value = 4025.5667
$ws.Range("N34").Value = -4429.566699999999

$ws.Range("H62").Value = 8697.571
$ws.Range("J62").Value = 7924
$ws.Range("L62").Value = 7924
$ws.Range("N62").Value = -9172

$ws.Range("H65").Value = 8697.571
$ws.Range("J65").Value = 7924
$ws.Range("L65").Value = 39620
$ws.Range("N65").Value = -45860

$ws.Range("H86").Value = 4597.6
$ws.Range("I86").Value = 4597.6
$ws.Range("K86").Value = 4597.6
$ws.Range("M86").Value = -3474.6

$ws.Range("H89").Value = 4597.6
$ws.Range("I89").Value = 4597.6
$ws.Range("K89").Value = 22988
$ws.Range("M89").Value = -17372

$ws.Range("H105").Value = 1544.3077
$ws.Range("I105").Value = 1564.75
$ws.Range("K105").Value = 1564.75
$ws.Range("M105").Value = 182.25

$ws.Range("H107").Value = 1018.73914
$ws.Range("I107").Value = 1100.25
$ws.Range("K107").Value = 1100.25
$ws.Range("M107").Value = 819.75

$ws.Range("H122").Value = 2411.077
$ws.Range("J122").Value = 6799.5
$ws.Range("L122").Value = 20398.5
$ws.Range("N122").Value = -25298.5

$ws.Range("H123").Value = 58000
$ws.Range("J123").Value = 58000
$ws.Range("L123").Value = 58000
$ws.Range("N123").Value = -67800

$ws.Range("H134").Value = 2606.389
$ws.Range("I134").Value = 1541.7
$ws.Range("K134").Value = 4625.1
$ws.Range("M134").Value = -2090.1

$ws.Range("H137").Value = 59740
$ws.Range("J137").Value = 59740
$ws.Range("L137").Value = 59740
$ws.Range("N137").Value = -69940

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 1363.24
$ws.Range("J107").Value = 1350.25
$ws.Range("L107").Value = 4050.75
$ws.Range("N107").Value = -7890.75

$ws.Range("H131").Value = 13571.667
$ws.Range("I131").Value = 1434
$ws.Range("J131").Value = 19640.5
$ws.Range("K131").Value = 4302
$ws.Range("L131").Value = 58921.5
$ws.Range("M131").Value = 738
$ws.Range("N131").Value = -69001.5

$ws.Range("H137").Value = 1714.579
$ws.Range("I137").Value = 904.6
$ws.Range("J137").Value = 2614.5557
$ws.Range("K137").Value = 2713.8
$ws.Range("L137").Value = 7843.6671
$ws.Range("M137").Value = 2386.2
$ws.Range("N137").Value = -18043.6671

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 7831.9614
$ws.Range("I126").Value = 9340.666999999999
$ws.Range("K126").Value = 28022.001
$ws.Range("M126").Value = -25552.001

$ws.Range("H132").Value = 181961.27
$ws.Range("I132").Value = 259200.25
$ws.Range("K132").Value = 777600.75
$ws.Range("M132").Value = -775070.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2226.1843
$ws.Range("I16").Value = 1987.8518
$ws.Range("J16").Value = 2811.182
$ws.Range("K16").Value = 1987.8518
$ws.Range("L16").Value = 2811.182
$ws.Range("M16").Value = -1817.8518
$ws.Range("N16").Value = -3151.182

$ws.Range("H59").Value = 50000
$ws.Range("J59").Value = 50000
$ws.Range("L59").Value = 50000
$ws.Range("N59").Value = -51308

$ws.Range("H93").Value = 1925.6666
$ws.Range("I93").Value = 1741.9286
$ws.Range("K93").Value = 1741.9286
$ws.Range("M93").Value = -493.9286

$ws.Range("H132").Value = 5641
$ws.Range("I132").Value = 3094.5715
$ws.Range("K132").Value = 9283.7145
$ws.Range("M132").Value = -6753.7145

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 9993
$ws.Range("J5").Value = 9993
$ws.Range("L5").Value = 9993
$ws.Range("N5").Value = -10217

$ws.Range("H81").Value = 11117611
$ws.Range("I81").Value = 4193.125
$ws.Range("J81").Value = 20008344
$ws.Range("K81").Value = 8386.25
$ws.Range("L81").Value = 40016688
$ws.Range("M81").Value = -7325.25
$ws.Range("N81").Value = -40018810

$ws.Range("H84").Value = 11117611
$ws.Range("I84").Value = 4193.125
$ws.Range("J84").Value = 20008344
$ws.Range("K84").Value = 41931.25
$ws.Range("L84").Value = 200083440
$ws.Range("M84").Value = -36627.25
$ws.Range("N84").Value = -200094048

$ws.Range("H126").Value = 1605.6666
$ws.Range("I126").Value = 1553.3077
$ws.Range("K126").Value = 4659.9231
$ws.Range("M126").Value = -2189.9231

$ws.Range("H136").Value = 3900.5789
$ws.Range("I136").Value = 2426
$ws.Range("J136").Value = 6428.4287
$ws.Range("K136").Value = 7278
$ws.Range("L136").Value = 19285.2861
$ws.Range("M136").Value = -4728
$ws.Range("N136").Value = -24385.2861

$ws.Range("H141").Value = 50000
$ws.Range("J141").Value = 50000
$ws.Range("L141").Value = 50000
$ws.Range("N141").Value = -60360
